$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlighted "group" rows: recolor fill from gray to yellow ---
# (rows 2-7 = group 1, rows 14-19 = group 3, both previously shaded gray)
$ws.Range("A2:D7").Interior.Color = 65535
$ws.Range("A14:D19").Interior.Color = 65535

# --- Mark two students in red font ---
$ws.Range("C7").Font.Color = 255
$ws.Range("C13").Font.Color = 255

# --- New column E: topic note for each group, placed on the group's first row ---
$ws.Range("E2").Value = "subgrid layout"
$ws.Range("E8").Value = "new viewport unit"
$ws.Range("E14").Value = "nesting layout"
$ws.Range("E20").Value = "video, audio"

# widen the new column
$ws.Columns("E").ColumnWidth = 30.83

# --- Move the active selection ---
[void]$ws.Range("G10").Select()

Write-Host "done"
